# Symlex VPN report_summary.xlsx edit:
# "created api symlex login brand,model module pos,neg test cases"
#
# A new worksheet "server_list" is inserted right before "get_support"
# (becoming the new active tab), and the "get_support" sheet's test
# summary numbers are updated to reflect one more failing test case.

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "server_list" sheet -----------------------
# Duplicate "get_support" (same D2:E8 layout, styles and merged title
# cell) and drop the copy in right before it; that's exactly where the
# new sheet belongs in the tab order.
$getSupport = $wb.Worksheets.Item("get_support")
$getSupport.Copy($getSupport)

$serverList = $wb.Worksheets.Item("get_support (2)")
$serverList.Name = "server_list"

# The original "get_support" worksheet object's position/index shifted
# once the copy was inserted in front of it, so re-resolve it by name
# instead of reusing the (now stale) $getSupport reference.
$getSupport = $wb.Worksheets.Item("get_support")

# --- Step 2: update "get_support" numbers first ------------------------
# (claims shared-string slot 44 = "TC_SYM_GSF_0021" ahead of the new
# sheet's own strings, matching the author's original edit order)
$getSupport.Range("E3").Value = 21
$getSupport.Range("E5").Value = 1
$getSupport.Range("E8").Value = "TC_SYM_GSF_0021"

# --- Step 3: fill in the new "server_list" sheet's own numbers --------
$serverList.Range("D2").Value = " Test Case Summary (07-03-24)"
$serverList.Range("E3").Value = 29
$serverList.Range("E4").Value = 25
$serverList.Range("E5").Value = 2
$serverList.Range("E6").Value = 2
$serverList.Range("E8").Value = "TC_SYM_SP_017,
TC_SYM_SP_029"

# --- Step 4: fix up selections / active tab ----------------------------
# get_support is no longer the selected tab, and its last selection
# moved from E6 to E8.
$getSupport.Activate()
$getSupport.Range("E8").Select()

# server_list becomes the active/selected sheet, with E4 selected.
$serverList.Activate()
$serverList.Range("E4").Select()
